$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The dataset gained a new weekly price observation for "Rabanito" at
# Vega Modelo de Temuco. This inserts a new row at position 53, pushing
# the existing rows 53-68 down to 54-69 (dimension grows to A1:R69),
# and fills the newly vacated row 53 with the new reading while the
# existing row 52 is updated to hold that new reading's date/volume
# (matching the author's target layout exactly).
$ws.Rows("53:53").Insert()

# Update row 52 (the anchor row) with the new observation's date & volume
$ws.Range("D52").Value = 44722
$ws.Range("J52").Value = 30

# Fill the newly-inserted row 53 with the data that used to sit in row 52
$ws.Range("A53").Value = 10
$ws.Range("B53").Value = "Vega Modelo de Temuco"
$ws.Range("C53").Value = "La Araucanía"
$ws.Range("D53").Value = 44623
$ws.Range("E53").Value = 9
$ws.Range("F53").Value = 300000001
$ws.Range("G53").Value = "Rabanito"
$ws.Range("H53").Value = "Sin especificar"
$ws.Range("I53").Value = "Primera"
$ws.Range("J53").Value = 50
$ws.Range("K53").Value = 7000
$ws.Range("L53").Value = 7000
$ws.Range("M53").Value = 7000
$ws.Range("N53").Value = "$/docena de paquetes"
$ws.Range("O53").Value = "Provincia de Cautín"
$ws.Range("P53").Value = 583
$ws.Range("Q53").Value = 12
$ws.Range("R53").Value = "Hortaliza"
